$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B: uniform change from 1.05 to 1.02 for data rows 2-25
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 2).Value = 1.02
}

# Columns C,D,E,F,I,J,K,L,M,N: new simulation result values per row
# Column indices: C=3 D=4 E=5 F=6 I=9 J=10 K=11 L=12 M=13 N=14
$colIndexes = @(3, 4, 5, 6, 9, 10, 11, 12, 13, 14)

$newValues = @{}
$newValues[2] = @(1.038832713560373, 1.046212595572148, 1.052225749829291, 1.059278117667195, 1.040228922949046, 1.043928010663119, 1.048978663420078, 1.054975070001183, 1.062008058181431, 1.018563602310288)
$newValues[3] = @(1.03974038314036, 1.046903100507909, 1.053077726659516, 1.060118304158643, 1.040405394124407, 1.044480845720832, 1.049481080248608, 1.055639749587014, 1.062662376942531, 1.018747886645176)
$newValues[4] = @(1.040328218689094, 1.047350022283729, 1.053629772242466, 1.060662506162755, 1.040518002300102, 1.044838439787375, 1.049805618780901, 1.056069955655988, 1.063085661189265, 1.018867052642491)
$newValues[5] = @(1.040575466161091, 1.047537935004309, 1.053862032390386, 1.060891417458649, 1.040564963609865, 1.044988740714488, 1.049941919742914, 1.056250840170799, 1.063263583583357, 1.018917130854643)
$newValues[6] = @(1.040616987180659, 1.047569487900482, 1.053901040408866, 1.060929860173182, 1.040572826362609, 1.045013975013318, 1.049964797318552, 1.056281212954072, 1.063293455951095, 1.018925538069971)
$newValues[7] = @(1.040331521946519, 1.047352533079466, 1.053632875008024, 1.060665564382896, 1.040518631290294, 1.044840448241432, 1.049807440573232, 1.056072372543777, 1.063088038702205, 1.018867721865891)
$newValues[8] = @(1.039139357806607, 1.046445929329755, 1.052513521106256, 1.05956194870131, 1.040288888977494, 1.044114869490045, 1.04914857246975, 1.05519967715291, 1.062229208870497, 1.018625898016108)
$newValues[9] = @(1.037042596657584, 1.044849361579161, 1.050546967640213, 1.057621488852503, 1.039871984245197, 1.042835384551658, 1.047983339284501, 1.053662810307522, 1.060715111041961, 1.018199192851902)
$newValues[10] = @(1.035647512262885, 1.043785742742214, 1.049239991320062, 1.056330810538029, 1.039585979891499, 1.041981843389702, 1.047203752480256, 1.052638942487593, 1.059705308257536, 1.017914357751486)
$newValues[11] = @(1.035044095122419, 1.043325384770063, 1.048675039551272, 1.055772656520118, 1.039460233006459, 1.041612134286761, 1.046865542660532, 1.052195781985176, 1.059267974688397, 1.017790939741102)
$newValues[12] = @(1.034820060083588, 1.043154418408681, 1.048465339828559, 1.055565442884151, 1.039413239544164, 1.041474790949826, 1.046739821014792, 1.05203120099983, 1.05910551859597, 1.017745084817424)
$newValues[13] = @(1.034868111771365, 1.04319108984053, 1.048510314383754, 1.055609885924514, 1.039423332707438, 1.041504252309029, 1.046766793034586, 1.052066502865781, 1.059140366451024, 1.017754921389569)
$newValues[14] = @(1.03502557424952, 1.043311251993668, 1.048657702667535, 1.055755525930244, 1.039456354331374, 1.041600781782178, 1.046855152410314, 1.05218217707917, 1.059254546231867, 1.017787149600757)
$newValues[15] = @(1.035122605494217, 1.043385292030351, 1.048748533182931, 1.055845274121215, 1.039476662238592, 1.041660254569037, 1.046909580925369, 1.052253451604513, 1.05932489475737, 1.017807004890496)
$newValues[16] = @(1.035687573159707, 1.043816299454299, 1.049277506045764, 1.056367868660293, 1.039594285223759, 1.042006377330891, 1.047226184897677, 1.05266835752071, 1.059734331034535, 1.017922546895385)
$newValues[17] = @(1.036042141001986, 1.044086712658848, 1.049609579641336, 1.05669587190484, 1.039667557493002, 1.042223459630556, 1.047424610768591, 1.052928666465784, 1.059991138571869, 1.017995001543065)
$newValues[18] = @(1.036249018369953, 1.044244458907787, 1.049803366765695, 1.056887259829961, 1.03971011210196, 1.042350068394647, 1.047540287042529, 1.053080517644051, 1.060140922026941, 1.018037255121926)
$newValues[19] = @(1.036319569014172, 1.044298249441297, 1.049869459097489, 1.056952529848747, 1.039724590914155, 1.042393236697564, 1.047579719078356, 1.05313229787727, 1.060191992907611, 1.018051661132132)
$newValues[20] = @(1.036004092604951, 1.044057697924551, 1.049573941530321, 1.056660673092672, 1.039659715084638, 1.04220016995296, 1.047403327973122, 1.052900735972527, 1.059963586374419, 1.017987228668611)
$newValues[21] = @(1.034979202677435, 1.043275866363236, 1.048614296378477, 1.055712635538783, 1.039446638160519, 1.041572356716086, 1.046829135401742, 1.052148113107201, 1.05922092342257, 1.01777765952315)
$newValues[22] = @(1.034335398415218, 1.042784478497508, 1.048011789479262, 1.055117201764474, 1.039311016968912, 1.041177528313662, 1.046467566048798, 1.051675075134251, 1.058753919203041, 1.017645825981718)
$newValues[23] = @(1.034676635392297, 1.043044954830981, 1.048331107710531, 1.055432791705626, 1.039383068599722, 1.041386843164423, 1.046659292687439, 1.051925825335214, 1.05900149242694, 1.017715719846391)
$newValues[24] = @(1.036021284840331, 1.044070808376144, 1.049590044567739, 1.056676577707203, 1.039663259299592, 1.042210693591943, 1.047412944944538, 1.052913356506635, 1.059976036051459, 1.017990740920734)
$newValues[25] = @(1.037584179523649, 1.045261985884645, 1.05105466056499, 1.058122629749449, 1.039981189556859, 1.043166264276866, 1.048285072759113, 1.054060008075679, 1.06110661919648, 1.018309572553644)

foreach ($r in $newValues.Keys) {
    $vals = $newValues[$r]
    for ($i = 0; $i -lt $colIndexes.Length; $i++) {
        $ws.Cells.Item($r, $colIndexes[$i]).Value = $vals[$i]
    }
}
